# Add two new rows (48 and 49) to the EOD worksheet, continuing the
# existing pattern of date/amount entries found in column A/B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValue = $ws.Range("A47").Text

$ws.Range("A48").Value = $dateValue
$ws.Range("B48").Value = 350

$ws.Range("A49").Value = $dateValue
$ws.Range("B49").Value = 350
